# The workbook records quarterly "A/B/C/D" sub-period rows for each year.
# This edit:
#   1) Swaps the "B" and "C" sub-period rows (columns A:E) within every
#      year block, so the row order becomes A, C, B, D.
#   2) Drops columns F ("...产销率") and G ("...销售量"), which duplicated
#      (derivable from) columns B and E, shrinking the used range from
#      A1:G69 to A1:E69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers holding the "B" sub-period and "C" sub-period for each of the
# 17 year-blocks (2000-2012, then 2016-2019), in order.
$rowsB = @(3, 7, 11, 15, 19, 23, 27, 31, 35, 39, 43, 47, 51, 55, 59, 63, 67)
$rowsC = @(4, 8, 12, 16, 20, 24, 28, 32, 36, 40, 44, 48, 52, 56, 60, 64, 68)

for ($i = 0; $i -lt $rowsB.Length; $i++) {
    $rB = $rowsB[$i]
    $rC = $rowsC[$i]

    $rngB = $ws.Range("A$rB`:E$rB")
    $rngC = $ws.Range("A$rC`:E$rC")

    $valB = $rngB.Value2
    $valC = $rngC.Value2

    $rngB.Value2 = $valC
    $rngC.Value2 = $valB
}

# Remove columns F and G entirely (shifts nothing left of F, just drops the
# two trailing columns and shrinks the sheet's used range / dimension).
$ws.Range("F:G").Delete()
